# Delete the row for code "0403 9099 39" (row 36) from Sheet1.
# This mirrors the author's edit: the "10 Digit Codes Introduced from 1
# January 2022" table no longer lists the 0403 9099 39 line, and all rows
# below it move up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(36).Select()
$ws.Rows.Item(36).Delete()
